# Updates cryptos list values per the upstream data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.081.57"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3
$ws.Range("D3").Value = "'2.551.76"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'575.89"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
$ws.Range("D6").Value = "'147.02"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -1.04%  "

# Row 9
$ws.Range("E9").Value = "  -1.14%  "

# Row 10
$ws.Range("D10").Value = "'5.51"
$ws.Range("E10").Value = "  -4.15%  "

# Row 11
$ws.Range("E11").Value = "  -0.57%  "

# Row 12
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  -0.67%  "

# Row 13
$ws.Range("D13").Value = "'27.19"
$ws.Range("E13").Value = "  -3.51%  "

# Row 14
$ws.Range("D14").Value = "'3.006.37"
$ws.Range("E14").Value = "  +0.26%  "

# Row 15
$ws.Range("D15").Value = "'62.982.71"
$ws.Range("E15").Value = "  -1.00%  "

# Row 16
$ws.Range("D16").Value = "'0.0000143"
$ws.Range("E16").Value = "  -0.74%  "

# Row 17
$ws.Range("D17").Value = "'2.544.69"
$ws.Range("E17").Value = "  +0.10%  "

# Row 18
$ws.Range("D18").Value = "'11.35"
$ws.Range("E18").Value = "  -1.47%  "

# Row 19
$ws.Range("D19").Value = "'336.20"
$ws.Range("E19").Value = "  -1.63%  "

# Row 20
$ws.Range("D20").Value = "'4.33"
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$ws.Range("D21").Value = "'6.75"
$ws.Range("E21").Value = "  -1.77%  "

# Row 23
$ws.Range("D23").Value = "'65.34"
$ws.Range("E23").Value = "  -1.05%  "

# Row 24
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("D25").Value = "'1.60"
$ws.Range("E25").Value = "  +1.61%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "'1.48"
$ws.Range("E27").Value = "  +4.56%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'8.34"
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = "  +3.35%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0813"
$ws.Range("E30").Value = "  -2.64%  "

# Row 31
$ws.Range("D31").Value = "'1.87"
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("D32").Value = "'178.20"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33
$ws.Range("D33").Value = "'1.55"
$ws.Range("E33").Value = "  -3.95%  "

# Row 34
$ws.Range("D34").Value = "'408.60"
$ws.Range("E34").Value = "  -3.49%  "

# Row 35
$ws.Range("D35").Value = "'19.10"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").Value = "'0.400"
$ws.Range("E36").Value = "  -1.55%  "

# Row 38
$ws.Range("D38").Value = "'4.34"
$ws.Range("E38").Value = "  -2.11%  "

# Row 39
$ws.Range("D39").Value = "'1.74"
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").Value = "'39.37"
$ws.Range("E41").Value = "  -3.26%  "

# Row 42
$ws.Range("D42").Value = "'151.26"
$ws.Range("E42").Value = "  -2.04%  "

# Row 43
$ws.Range("D43").Value = "'3.75"
$ws.Range("E43").Value = "  -0.96%  "

# Row 44
$ws.Range("D44").Value = "'20.85"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'0.0538"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46
$ws.Range("D46").Value = "'0.603"
$ws.Range("E46").Value = "  -0.99%  "

# Row 47
$ws.Range("D47").Value = "'0.0964"
$ws.Range("E47").Value = "  -0.30%  "

# Row 48
$ws.Range("E48").Value = "  +2.89%  "

# Row 49
$ws.Range("D49").Value = "'18.22"
$ws.Range("E49").Value = "  -2.62%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'1.71"
$ws.Range("E50").Value = "  -7.94%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.28"
$ws.Range("E51").Value = "  +0.18%  "

Write-Output "Updated cryptos list values."